# 20241108 - Data Science Personal Log.xlsx
# "updated for Nov 20."
#
# 1. Hockey: fill in rows 31 & 32 (games on Nov 19 / Nov 20 2024)
# 2. New "Job Aps" sheet inserted between Bowling and Dates with a job-application tracker
# 3. Poutine: add a new tasting entry (row 20)
# 4. Books: log a new book started (row 11)
# 5. Selection / active-sheet bookkeeping (Books becomes the active tab)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Hockey sheet - rows 31 & 32
# ---------------------------------------------------------------------------
$hockey = $wb.Worksheets.Item("Hockey")

# Row 31 (A31/B31 already populated: 2024-11-19, George Bell)
$hockey.Cells.Item(31, 3).Value = 7.9
$hockey.Cells.Item(31, 4).Value = 0
$hockey.Cells.Item(31, 5).Value = 1
$hockey.Cells.Item(31, 6).Value = 1
$hockey.Cells.Item(31, 7).Value = 2
$hockey.Cells.Item(31, 8).Value = 3
$hockey.Cells.Item(31, 9).Value = "played well. Good defense and passes. Decent shot"
$hockey.Cells.Item(31, 10).Value = "offense not the best"

# Row 32 (A32 already populated: 2024-11-20)
$hockey.Cells.Item(32, 2).Value = "ASHL"
$hockey.Cells.Item(32, 3).Value = 8.4
$hockey.Cells.Item(32, 4).Value = 0
$hockey.Cells.Item(32, 5).Value = 0
$hockey.Cells.Item(32, 6).Value = 1
$hockey.Cells.Item(32, 7).Value = 1
$hockey.Cells.Item(32, 8).Value = 3
$hockey.Cells.Item(32, 9).Value = "passing was really good. Maybe passed 90% conversion. Won against a bunch of old guys"
$hockey.Cells.Item(32, 10).Value = "fell a couple times"

$hockey.Range("J30").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. New "Job Aps" sheet, positioned right after Bowling (before Dates)
# ---------------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("Bowling")
$jobAps = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $bowling)
$jobAps.Name = "Job Aps"

# Header row
$jobAps.Cells.Item(3, 2).Value = "Date"
$jobAps.Cells.Item(3, 3).Value = "Job Title"
$jobAps.Cells.Item(3, 4).Value = "Company"
$jobAps.Cells.Item(3, 5).Value = "End Date"
$jobAps.Cells.Item(3, 6).Value = "Other Info"
$jobAps.Cells.Item(3, 7).Value = "Info"
$jobAps.Cells.Item(3, 9).Value = "Link"

# Row 4
$jobAps.Cells.Item(4, 2).Value = [datetime]"2018-11-19"
$jobAps.Cells.Item(4, 3).Value = "Analyst, Cancer Screening Analytics"
$jobAps.Cells.Item(4, 4).Value = "Ontario Health"
$jobAps.Cells.Item(4, 5).Value = "November 23, 2024 (4 days left to apply)"
$jobAps.Cells.Item(4, 6).Value = "R107915"
$jobAps.Cells.Item(4, 8).Value = "Direct Website"
$jobAps.Cells.Item(4, 9).Value = "https://oh.wd3.myworkdayjobs.com/en-US/OH/job/Toronto-ON/Analyst--Health-Data-II_R107955"
$jobAps.Cells.Item(4, 10).Value = "-"

# Row 5
$jobAps.Cells.Item(5, 2).Value = [datetime]"2018-11-19"
$jobAps.Cells.Item(5, 3).Value = "Analyst, Health Data"
$jobAps.Cells.Item(5, 4).Value = "Ontario Health"
$jobAps.Cells.Item(5, 5).Value = "November 26, 2024 (7 days left to apply)"
$jobAps.Cells.Item(5, 8).Value = "Direct Website"
$jobAps.Cells.Item(5, 9).Value = "https://oh.wd3.myworkdayjobs.com/en-US/OH/job/Toronto-ON/Analyst--Health-Data-II_R107955"
$jobAps.Cells.Item(5, 10).Value = "-"

# Row 6
$jobAps.Cells.Item(6, 2).Value = [datetime]"2018-11-20"
$jobAps.Cells.Item(6, 3).Value = "Business Analyst"
$jobAps.Cells.Item(6, 4).Value = "Atlantis IT group"
$jobAps.Cells.Item(6, 8).Value = "Indeed"
$jobAps.Cells.Item(6, 9).Value = "https://ca.indeed.com/?vjk=dd81739d3c279da3"
$jobAps.Cells.Item(6, 10).Value = "-"

# Row 7
$jobAps.Cells.Item(7, 2).Value = [datetime]"2018-11-20"
$jobAps.Cells.Item(7, 3).Value = "Financial Planning Analyst - job post"
$jobAps.Cells.Item(7, 4).Value = "Spectrum Health Care"
$jobAps.Cells.Item(7, 8).Value = "Indeed"
$jobAps.Cells.Item(7, 9).Value = "https://ca.indeed.com/jobs?q=data&l=Toronto%2C+ON&vjk=14389fc4c61f4bb2"
$jobAps.Cells.Item(7, 10).Value = "-"

# Row 8
$jobAps.Cells.Item(8, 2).Value = [datetime]"2018-11-20"
$jobAps.Cells.Item(8, 3).Value = "Economist / Research Analyst"
$jobAps.Cells.Item(8, 4).Value = "UFCW Canada"
$jobAps.Cells.Item(8, 8).Value = "Indeed"
$jobAps.Cells.Item(8, 9).Value = "https://ca.indeed.com/jobs?q=data&l=Toronto%2C+ON&vjk=14389fc4c61f4bb2"
$jobAps.Cells.Item(8, 10).Value = "-"

$jobAps.Columns.Item(2).ColumnWidth = 10.42578125

$jobAps.Range("J15").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Poutine sheet - row 20 (new tasting)
# ---------------------------------------------------------------------------
$poutine = $wb.Worksheets.Item("Poutine")
$poutine.Cells.Item(20, 2).Value = [datetime]"2024-11-20"
$poutine.Cells.Item(20, 3).Value = "Burger Legend Etobicoke"
$poutine.Cells.Item(20, 4).Value = 8.9
$poutine.Cells.Item(20, 5).Value = 13
$poutine.Cells.Item(20, 6).Value = "M"
$poutine.Cells.Item(20, 7).Value = "8 30pm"
$poutine.Cells.Item(20, 8).Value = "Regular"
$poutine.Cells.Item(20, 9).Value = "Surpringly good - impressive. Gravy, everything, good"

$poutine.Range("B21").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Books sheet - row 11 (new book started)
# ---------------------------------------------------------------------------
$books = $wb.Worksheets.Item("Books")
$books.Cells.Item(11, 4).Value = "Technofeudalism"

# Books becomes the active tab, selection lands on E11
$books.Activate() | Out-Null
$books.Range("E11").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. Golf sheet - selection only, no data change
# ---------------------------------------------------------------------------
$golf = $wb.Worksheets.Item("Golf")
$golf.Range("I17").Select() | Out-Null

# Leave Books as the final active sheet/tab (matches activeTab in the target file)
$books.Activate() | Out-Null
